$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.264243666666667
$ws.Range("H2").Value = 6.792731
$ws.Range("I2").Value = 0.4114976873616865
$ws.Range("J2").Value = 0.4114976873616865
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.715641666666667
$ws.Range("N2").Value = 8.146925
$ws.Range("O2").Value = 0.1838906555826552
$ws.Range("P2").Value = 0.1838906555826552
$ws.Range("Q2").Value = 6.148874444686111
$ws.Range("R2").Value = 55.339870002175
$ws.Range("S2").Value = 0.07567057949968704
$ws.Range("T2").Value = 0.07567057949968703
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.264243666666667
$ws.Range("H3").Value = 6.792731
$ws.Range("I3").Value = 0.4114976873616865
$ws.Range("J3").Value = 0.4114976873616865
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 9.604052333333334
$ws.Range("N3").Value = 28.812157
$ws.Range("O3").Value = 0.65034187002831
$ws.Range("P3").Value = 0.6503418700283099
$ws.Range("Q3").Value = 21.74591467008522
$ws.Range("R3").Value = 195.713232030767
$ws.Range("S3").Value = 0.2676141755111241
$ws.Range("T3").Value = 0.267614175511124
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.264243666666667
$ws.Range("H4").Value = 6.792731
$ws.Range("I4").Value = 0.4114976873616865
$ws.Range("J4").Value = 0.4114976873616865
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.025443666666667
$ws.Range("N4").Value = 3.076331
$ws.Range("O4").Value = 0.06943828798954764
$ws.Range("P4").Value = 0.06943828798954764
$ws.Range("Q4").Value = 2.321854327773444
$ws.Range("R4").Value = 20.896688949961
$ws.Range("S4").Value = 0.02857369492205362
$ws.Range("T4").Value = 0.02857369492205362
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.264243666666667
$ws.Range("H5").Value = 6.792731
$ws.Range("I5").Value = 0.4114976873616865
$ws.Range("J5").Value = 0.4114976873616865
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.422560333333333
$ws.Range("N5").Value = 4.267681000000001
$ws.Range("O5").Value = 0.09632918639948715
$ws.Range("P5").Value = 0.09632918639948714
$ws.Range("Q5").Value = 3.221023225201223
$ws.Range("R5").Value = 28.989209026811
$ws.Range("S5").Value = 0.03963923742882178
$ws.Range("T5").Value = 0.03963923742882178
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.089228666666667
$ws.Range("H6").Value = 6.267686
$ws.Range("I6").Value = 0.3796909216792509
$ws.Range("J6").Value = 0.3796909216792509
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.715641666666667
$ws.Range("N6").Value = 8.146925
$ws.Range("O6").Value = 0.1838906555826552
$ws.Range("P6").Value = 0.1838906555826552
$ws.Range("Q6").Value = 5.673596418394444
$ws.Range("R6").Value = 51.06236776555
$ws.Range("S6").Value = 0.06982161250638005
$ws.Range("T6").Value = 0.06982161250638004
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.089228666666667
$ws.Range("H7").Value = 6.267686
$ws.Range("I7").Value = 0.3796909216792509
$ws.Range("J7").Value = 0.3796909216792509
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 9.604052333333334
$ws.Range("N7").Value = 28.812157
$ws.Range("O7").Value = 0.65034187002831
$ws.Range("P7").Value = 0.6503418700283099
$ws.Range("Q7").Value = 20.06506145096689
$ws.Range("R7").Value = 180.585553058702
$ws.Range("S7").Value = 0.2469289040376566
$ws.Range("T7").Value = 0.2469289040376566
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.089228666666667
$ws.Range("H8").Value = 6.267686
$ws.Range("I8").Value = 0.3796909216792509
$ws.Range("J8").Value = 0.3796909216792509
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.025443666666667
$ws.Range("N8").Value = 3.076331
$ws.Range("O8").Value = 0.06943828798954764
$ws.Range("P8").Value = 0.06943828798954764
$ws.Range("Q8").Value = 2.142386304451778
$ws.Range("R8").Value = 19.281476740066
$ws.Range("S8").Value = 0.0263650875665806
$ws.Range("T8").Value = 0.0263650875665806
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.089228666666667
$ws.Range("H9").Value = 6.267686
$ws.Range("I9").Value = 0.3796909216792509
$ws.Range("J9").Value = 0.3796909216792509
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.422560333333333
$ws.Range("N9").Value = 4.267681000000001
$ws.Range("O9").Value = 0.09632918639948715
$ws.Range("P9").Value = 0.09632918639948714
$ws.Range("Q9").Value = 2.972053828462889
$ws.Range("R9").Value = 26.748484456166
$ws.Range("S9").Value = 0.03657531756863364
$ws.Range("T9").Value = 0.03657531756863363
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.1905406666666667
$ws.Range("H10").Value = 0.571622
$ws.Range("I10").Value = 0.03462835949856721
$ws.Range("J10").Value = 0.03462835949856721
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.715641666666667
$ws.Range("N10").Value = 8.146925
$ws.Range("O10").Value = 0.1838906555826552
$ws.Range("P10").Value = 0.1838906555826552
$ws.Range("Q10").Value = 0.5174401735944444
$ws.Range("R10").Value = 4.656961562349999
$ws.Range("S10").Value = 0.006367831729943391
$ws.Range("T10").Value = 0.00636783172994339
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.1905406666666667
$ws.Range("H11").Value = 0.571622
$ws.Range("I11").Value = 0.03462835949856721
$ws.Range("J11").Value = 0.03462835949856721
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 9.604052333333334
$ws.Range("N11").Value = 28.812157
$ws.Range("O11").Value = 0.65034187002831
$ws.Range("P11").Value = 0.6503418700283099
$ws.Range("Q11").Value = 1.829962534294889
$ws.Range("R11").Value = 16.469662808654
$ws.Range("S11").Value = 0.02252027207231079
$ws.Range("T11").Value = 0.02252027207231079
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.1905406666666667
$ws.Range("H12").Value = 0.571622
$ws.Range("I12").Value = 0.03462835949856721
$ws.Range("J12").Value = 0.03462835949856721
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1.025443666666667
$ws.Range("N12").Value = 3.076331
$ws.Range("O12").Value = 0.06943828798954764
$ws.Range("P12").Value = 0.06943828798954764
$ws.Range("Q12").Value = 0.1953887198757778
$ws.Range("R12").Value = 1.758498478882
$ws.Range("S12").Value = 0.002404533999467098
$ws.Range("T12").Value = 0.002404533999467098
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.1905406666666667
$ws.Range("H13").Value = 0.571622
$ws.Range("I13").Value = 0.03462835949856721
$ws.Range("J13").Value = 0.03462835949856721
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.422560333333333
$ws.Range("N13").Value = 4.267681000000001
$ws.Range("O13").Value = 0.09632918639948715
$ws.Range("P13").Value = 0.09632918639948714
$ws.Range("Q13").Value = 0.2710555942868889
$ws.Range("R13").Value = 2.439500348582
$ws.Range("S13").Value = 0.003335721696845932
$ws.Range("T13").Value = 0.003335721696845932
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.9584326666666668
$ws.Range("H14").Value = 2.875298
$ws.Range("I14").Value = 0.1741830314604954
$ws.Range("J14").Value = 0.1741830314604954
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 2.715641666666667
$ws.Range("N14").Value = 8.146925
$ws.Range("O14").Value = 0.1838906555826552
$ws.Range("P14").Value = 0.1838906555826552
$ws.Range("Q14").Value = 2.602759684294445
$ws.Range("R14").Value = 23.42483715865
$ws.Range("S14").Value = 0.03203063184664477
$ws.Range("T14").Value = 0.03203063184664477
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.9584326666666668
$ws.Range("H15").Value = 2.875298
$ws.Range("I15").Value = 0.1741830314604954
$ws.Range("J15").Value = 0.1741830314604954
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 9.604052333333334
$ws.Range("N15").Value = 28.812157
$ws.Range("O15").Value = 0.65034187002831
$ws.Range("P15").Value = 0.6503418700283099
$ws.Range("Q15").Value = 9.204837488642891
$ws.Range("R15").Value = 82.84353739778601
$ws.Range("S15").Value = 0.1132785184072186
$ws.Range("T15").Value = 0.1132785184072185
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.9584326666666668
$ws.Range("H16").Value = 2.875298
$ws.Range("I16").Value = 0.1741830314604954
$ws.Range("J16").Value = 0.1741830314604954
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1.025443666666667
$ws.Range("N16").Value = 3.076331
$ws.Range("O16").Value = 0.06943828798954764
$ws.Range("P16").Value = 0.06943828798954764
$ws.Range("Q16").Value = 0.9828187079597779
$ws.Range("R16").Value = 8.845368371638001
$ws.Range("S16").Value = 0.01209497150144632
$ws.Range("T16").Value = 0.01209497150144632
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.9584326666666668
$ws.Range("H17").Value = 2.875298
$ws.Range("I17").Value = 0.1741830314604954
$ws.Range("J17").Value = 0.1741830314604954
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 1.422560333333333
$ws.Range("N17").Value = 4.267681000000001
$ws.Range("O17").Value = 0.09632918639948715
$ws.Range("P17").Value = 0.09632918639948714
$ws.Range("Q17").Value = 1.363428293770889
$ws.Range("R17").Value = 12.270854643938
$ws.Range("S17").Value = 0.0167789097051858
$ws.Range("T17").Value = 0.0167789097051858
